# AutoCommit_12 декабря 2023 г. 11:39:49_SibNout2023
# Fills in homework grades (5, occasionally 4) for several students and
# records a missed-exam variant for one student, then moves the active
# selection from S8 to S9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Cells whose fill/format must change to the "graded" look (style 10,
#    the light/white filled look already used for Дз8/Дз9 elsewhere).
#    We copy the format from J8, which already carries that exact style,
#    onto each destination range, then write the grade values.
# ---------------------------------------------------------------------
$styleSrcFilled = $ws.Range("J8")
$styleSrcFilled.Copy()
$ws.Range("J10:K10").PasteSpecial($xlPasteFormats)
$ws.Range("H12:K12").PasteSpecial($xlPasteFormats)
$ws.Range("J13:K13").PasteSpecial($xlPasteFormats)
$ws.Range("J14:K14").PasteSpecial($xlPasteFormats)
$ws.Range("J15:K15").PasteSpecial($xlPasteFormats)
$ws.Range("I20").PasteSpecial($xlPasteFormats)
$ws.Range("J26:K26").PasteSpecial($xlPasteFormats)
$ws.Range("C30:K30").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Style 3 (center/wrap look used for Дз7 in some rows) for I15.
$ws.Range("I17").Copy()
$ws.Range("I15").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Style 11 (filled look, no wrap) for L26.
$ws.Range("I19").Copy()
$ws.Range("L26").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Grade values. Cells above already have their format set; the rest
#    keep their existing formatting untouched.
# ---------------------------------------------------------------------

# Балашова Алиса (row 9) - fills in Дз2..Дз9, last one is a 4.
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 4

# Row 10
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 5

# Row 12
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 5

# Row 13
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("R13").Value = 5

# Row 14
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 5
$ws.Range("R14").Value = 5

# Row 15
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 5

# Row 20
$ws.Range("I20").Value = 5

# Row 25
$ws.Range("R25").Value = 5

# Row 26
$ws.Range("G26").Value = 5
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 5
$ws.Range("R26").Value = 5

# Row 29 - Вариант for a student that had missed it.
$ws.Range("L29").Value = 1

# Row 30 - Саргас Кирилл, fills in Дз1.1..Дз9.
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 5
$ws.Range("R30").Value = 4

# ---------------------------------------------------------------------
# 3) Move the active selection from S8 to S9.
# ---------------------------------------------------------------------
$ws.Range("S9").Select()
